$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the cryptos list.
# NumberFormat is forced to Text ("@") before assigning so Excel does
# not reinterpret values like "20.515.17" or "4.920" as numbers/dates
# and strip meaningful formatting (trailing zeros, grouping dots, etc).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.515.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.466.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9568"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.93"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3601"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3073"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.066"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06619"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.463"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.122"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9572"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001019"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.465.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05940"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.458"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.15"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.15%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.538.85"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.098"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.11"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.625.84"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.908"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07978"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.920"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7998"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.222"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.467"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05758"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.693"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02043"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9575"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.32"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1860"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.308"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5245"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.514"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.11"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5183"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.798"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06436"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9835"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.43%  "
